$d = $word.ActiveDocument

$oldSnippet = "Kampagnendaten 2018"
$newText = "Kampagnendaten Leo: 14. bis 23. April, 14. bis 23. Mai"
$anchor = " 30. Oktober"

# Collect indices (1-based) of paragraphs that need editing first, since
# modifying the document while iterating live collections can be unreliable.
$targets = New-Object System.Collections.ArrayList
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx++
    $t = $p.Range.Text
    if ($t -like "*$oldSnippet*") {
        [void]$targets.Add($idx)
    }
}

# Process from the last paragraph to the first so earlier indices/offsets
# remain valid while we edit later ones.
for ($k = $targets.Count - 1; $k -ge 0; $k--) {
    $pIndex = $targets[$k]
    $p = $d.Paragraphs.Item($pIndex)
    $rng = $p.Range
    $pStart = $rng.Start
    $pEnd = $rng.End - 1   # exclude the paragraph mark

    $fullText = $rng.Text
    $spaceIdx = $fullText.IndexOf($anchor)
    $spaceStart = $pStart + $spaceIdx
    $spaceEnd = $spaceStart + 1

    # Delete the trailing part (date range text) after the lone space run.
    $afterRng = $d.Range($spaceEnd, $pEnd)
    if ($afterRng.Start -lt $afterRng.End) {
        $afterRng.Delete()
    }

    # Delete the leading part (title + constellation name) before the lone
    # space run; this run carries no explicit run formatting (w:rPr), which
    # is what we want to keep for the replacement text.
    $beforeRng = $d.Range($pStart, $spaceStart)
    if ($beforeRng.Start -lt $beforeRng.End) {
        $beforeRng.Delete()
    }

    # Now the paragraph contains just the single, unformatted space run.
    # Replace its text with the new campaign-dates string.
    $p2 = $d.Paragraphs.Item($pIndex)
    $onlyRng = $d.Range($p2.Range.Start, $p2.Range.End - 1)
    $onlyRng.Text = $newText
}
